$d = $word.ActiveDocument

$d.Content.Find.Execute("-2.4954", $true, $false, $false, $false, $false, $true, 1, $false, "-2.4883", 2)
$d.Content.Find.Execute("0.0126", $true, $false, $false, $false, $false, $true, 1, $false, "0.0128", 2)
$d.Content.Find.Execute("-1.3481", $true, $false, $false, $false, $false, $true, 1, $false, "-1.3382", 2)
$d.Content.Find.Execute("0.1776", $true, $false, $false, $false, $false, $true, 1, $false, "0.1808", 2)
$d.Content.Find.Execute("10.0920", $true, $false, $false, $false, $false, $true, 1, $false, "10.0974", 2)
